# Updated cryptos list on Thu Apr 11 06:09:32 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# coinranking.com crypto table, and swaps the XRP / LidoStakedEther rows
# (rows 7 and 8) to reflect their new rank order.
#
# Cells are written as literal text (not auto-coerced numbers) to match the
# source data, which stores every Price/Volume cell as a string. A
# Formula + Copy/PasteSpecial(xlPasteValues) round-trip is used instead of a
# plain .Value assignment so Excel's "looks like a number" auto-conversion
# never kicks in and no incidental number-format style gets attached to the
# cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text {
    param($row, $col, [string]$value)
    $cell = $ws.Cells.Item($row, $col)
    $escaped = $value -replace '"', '""'
    $cell.Formula = "=""$escaped"""
    $cell.Copy()
    $cell.PasteSpecial(-4163)   # xlPasteValues
}

function Set-Pct {
    param($row, $value)
    Set-Text $row 5 "  $value  "
}

# Row 2 - Bitcoin
Set-Text 2 4 "70.684.83"
Set-Pct  2 "+1.99%"

# Row 3 - Ethereum
Set-Text 3 4 "3.567.99"
Set-Pct  3 "+0.95%"

# Row 4 - TetherUSD
Set-Text 4 4 "0.999"
Set-Pct  4 "-0.09%"

# Row 5 - BNB
Set-Text 5 4 "610.80"
Set-Pct  5 "+5.01%"

# Row 6 - Solana
Set-Text 6 4 "174.45"
Set-Pct  6 "+1.39%"

# Rows 7/8 - XRP and LidoStakedEther swap positions
Set-Text 7 2 "LidoStakedEther"
Set-Text 7 3 "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
Set-Text 7 4 "3.565.87"
Set-Pct  7 "+1.07%"

Set-Text 8 2 "XRP"
Set-Text 8 3 "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-Text 8 4 "0.616"
Set-Pct  8 "+1.19%"

# Row 10 - Dogecoin
Set-Text 10 4 "0.197"
Set-Pct  10 "+3.81%"

# Row 11 - Toncoin
Set-Text 11 4 "7.61"
Set-Pct  11 "+13.60%"

# Row 12 - Cardano
Set-Text 12 4 "0.587"
Set-Pct  12 "+0.26%"

# Row 13 - Avalanche
Set-Text 13 4 "46.85"
Set-Pct  13 "-1.25%"

# Row 14 - ShibaInu
Set-Text 14 4 "0.0000278"
Set-Pct  14 "+1.04%"

# Row 15 - WrappedliquidstakedEther2.0
Set-Text 15 4 "4.147.87"
Set-Pct  15 "+1.37%"

# Row 16 - Polkadot
Set-Text 16 4 "8.43"
Set-Pct  16 "-1.45%"

# Row 17 - BitcoinCash
Set-Text 17 4 "616.06"
Set-Pct  17 "-1.78%"

# Row 18 - WrappedEther
Set-Text 18 4 "3.563.19"
Set-Pct  18 "+0.81%"

# Row 19 - WrappedBTC
Set-Text 19 4 "70.754.84"
Set-Pct  19 "+2.13%"

# Row 20 - TRON (price unchanged)
Set-Pct 20 "-2.24%"

# Row 21 - Chainlink
Set-Text 21 4 "17.43"
Set-Pct  21 "-0.51%"

# Row 22 - Polygon
Set-Text 22 4 "0.889"
Set-Pct  22 "-0.06%"

# Row 23 - Uniswap
Set-Text 23 4 "9.44"
Set-Pct  23 "-15.81%"

# Row 24 - InternetComputer(DFINITY)
Set-Text 24 4 "16.06"
Set-Pct  24 "+0.57%"

# Row 25 - Litecoin
Set-Text 25 4 "97.38"
Set-Pct  25 "-0.32%"

# Row 26 - PancakeSwap (price unchanged)
Set-Pct 26 "+0.32%"

# Row 27 - Dai (price unchanged)
Set-Pct 27 "+0.04%"

# Row 28 - ImmutableX
Set-Text 28 4 "2.64"
Set-Pct  28 "+0.07%"

# Row 29 - EthereumClassic
Set-Text 29 4 "33.54"
Set-Pct  29 "+2.02%"

# Row 30 - RenderToken
Set-Text 30 4 "9.14"
Set-Pct  30 "-1.75%"

# Row 31 - Filecoin
Set-Text 31 4 "8.54"
Set-Pct  31 "-0.26%"

# Row 32 - Stacks
Set-Text 32 4 "3.06"
Set-Pct  32 "-3.07%"

# Row 33 - NEARProtocol
Set-Text 33 4 "7.01"
Set-Pct  33 "+0.20%"

# Row 34 - Mantle (price unchanged)
Set-Pct 34 "-2.14%"

# Row 35 - Bittensor
Set-Text 35 4 "598.54"
Set-Pct  35 "-5.50%"

# Row 36 - dogwifhat
Set-Text 36 4 "3.73"
Set-Pct  36 "+6.63%"

# Row 37 - Hedera
Set-Text 37 4 "0.102"
Set-Pct  37 "-1.04%"

# Row 38 - Cosmos
Set-Text 38 4 "10.85"
Set-Pct  38 "+0.63%"

# Row 39 - VeChain
Set-Text 39 4 "0.0482"
Set-Pct  39 "+7.00%"

# Row 40 - OKB
Set-Text 40 4 "57.42"
Set-Pct  40 "-0.03%"

# Row 41 - FirstDigitalUSD
Set-Text 41 4 "0.999"
Set-Pct  41 "+0.01%"

# Row 42 - Kaspa (price unchanged)
Set-Pct 42 "+4.02%"

# Row 43 - Maker
Set-Text 43 4 "3.392.75"
Set-Pct  43 "+0.00%"

# Row 44 - TheGraph
Set-Text 44 4 "0.322"
Set-Pct  44 "-2.25%"

# Row 45 - InjectiveProtocol
Set-Text 45 4 "33.37"
Set-Pct  45 "+1.33%"

# Row 46 - ThetaToken
Set-Text 46 4 "3.00"
Set-Pct  46 "+7.99%"

# Row 47 - PEPE
Set-Text 47 4 "0.0₃0709"
Set-Pct  47 "+1.16%"

# Row 48 - Fetch.AI
Set-Text 48 4 "2.62"
Set-Pct  48 "+2.16%"

# Row 49 - Stellar (price unchanged)
Set-Pct 49 "+0.12%"

# Row 50 - Monero
Set-Text 50 4 "133.74"
Set-Pct  50 "+1.23%"
